# Apply cryptocurrency price/volume updates per commit on Thu Nov  2 19:52:41 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Text)
    $Worksheet.Range($Address).NumberFormat = "@"
    $Worksheet.Range($Address).Value = $Text
    $Worksheet.Range($Address).Style = "Normal"
}

$ws.Range("D2").Value = '35.100.88'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.818.09'
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  +0.15%  '
Set-TextValue $ws "D5" '233.28'
Set-TextValue $ws "D6" '0.614'
$ws.Range("E6").Value = '  +0.61%  '
Set-TextValue $ws "D7" '1.00'
$ws.Range("E7").Value = '  +0.12%  '
Set-TextValue $ws "D8" '41.23'
$ws.Range("E8").Value = '  -2.87%  '
Set-TextValue $ws "D9" '0.325'
$ws.Range("E9").Value = '  +7.63%  '
Set-TextValue $ws "D10" '0.0686'
$ws.Range("E10").Value = '  +0.06%  '
Set-TextValue $ws "D11" '0.1000'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '2.081.91'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").Value = '1.825.87'
$ws.Range("E13").Value = '  -0.61%  '
Set-TextValue $ws "D14" '11.10'
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws "D15" '0.661'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws "D16" '4.68'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '35.055.68'
$ws.Range("E17").Value = '  +1.03%  '
Set-TextValue $ws "D18" '69.65'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").Value = '0.0₃0792'
$ws.Range("E19").Value = '  +0.97%  '
Set-TextValue $ws "D20" '239.85'
$ws.Range("E20").Value = '  -1.38%  '
Set-TextValue $ws "D21" '11.88'
$ws.Range("E21").Value = '  -1.77%  '
Set-TextValue $ws "D22" '4.67'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("E24").Value = '  +3.09%  '
Set-TextValue $ws "D25" '172.93'
$ws.Range("E25").Value = '  +0.49%  '
Set-TextValue $ws "D26" '7.87'
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("E29").Value = '  +25.82%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("B31").Value = 'EURNeutrino'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn'
$ws.Range("D31").Value = '3.336.26'
$ws.Range("E31").Value = '  +37.31%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D32" '4.06'
$ws.Range("E32").Value = '  +4.00%  '
$ws.Range("E33").Value = '  +5.47%  '
Set-TextValue $ws "D34" '3.98'
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E35").Value = '  -4.57%  '
$ws.Range("E36").Value = '  +7.55%  '
Set-TextValue $ws "D37" '93.13'
$ws.Range("E37").Value = '  +3.59%  '
$ws.Range("E38").Value = '  +3.41%  '
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").Value = '1.310.85'
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("E42").Value = '  +0.46%  '
Set-TextValue $ws "D43" '14.76'
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("E44").Value = '  -4.54%  '
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("E46").Value = '  -2.09%  '
Set-TextValue $ws "D47" '6.36'
$ws.Range("E47").Value = '  +5.96%  '
Set-TextValue $ws "D48" '0.0512'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").Value = '1.997.16'
$ws.Range("E49").Value = '  -0.49%  '
Set-TextValue $ws "D50" '1.00'
$ws.Range("E50").Value = '  +0.16%  '
Set-TextValue $ws "D51" '0.0652'
$ws.Range("E51").Value = '  +6.50%  '
